# Test data change: add a "user_mobile" column (used for the "revoke license"
# test case) with the same mobile number repeated for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column V, row 1.
$ws.Range("V1").Value = "user_mobile"

# New mobile-number values for the three existing data rows (2-4).
$ws.Range("V2").Value = 4741854178
$ws.Range("V3").Value = 4741854178
$ws.Range("V4").Value = 4741854178

# The new column uses wrap text like its neighbouring data columns.
$ws.Range("V1:V4").WrapText = $true

# Try to reproduce the saved view state (best effort - not all view state is
# round-tripped by this host, but the calls are harmless if ignored).
$ws.Range("O1").Select()
$excel.ActiveWindow.ScrollColumn = 15
$ws.Range("V7").Select()
